$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Column B keeps "Algorithm"; columns C..L now hold mean/std pairs per horizon.
$headers = @(
    "Algorithm",
    "One Year Alt mean",
    "One Year Alt std",
    "Two Year Alt mean",
    "Two Year Alt std",
    "Three Year Alt mean",
    "Three Year Alt std",
    "Five Year Alt mean",
    "Five Year Alt std",
    "Ten Year Alt mean",
    "Ten Year Alt std"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    # headers start at column B (2) through column L (12)
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Columns H:L are brand new; copy the header formatting (bold, centered, bordered)
# from an existing header cell (B1) so the new cells render the same way.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows ---
# Each entry: algorithm name, then mean/std pairs for 1/2/3/5/10 year horizons (columns C..L)
$data = @(
    @("LR",    0.8269202200162644, 0.016265879147472,  0.8057980858379701, 0.01477018580232867, 0.789381124196361,  0.02072751941951771, 0.7722691340215022, 0.03297623102044418, 0.7530720308656191, 0.02301328700677731),
    @("LDA",   0.832217331559589,  0.01521574122255964, 0.8133104613414609, 0.01525340928887356, 0.801842049379341,  0.01942914514144753, 0.7878715678587421, 0.03593378303954588, 0.7711447119078412, 0.02557462646890845),
    @("KNN",   0.7661770929576098, 0.01016923241557916, 0.7704784661953039, 0.01400620361017207, 0.7679255817092716, 0.02562291102390764, 0.7776337348231441, 0.02379357330258653, 0.7751859513499821, 0.01907350294624797),
    @("DTREE", 0.7651556527242598, 0.02511881756105357, 0.7541573848414977, 0.02759614782270312, 0.7490592553412048, 0.02584965936616514, 0.7565132140232967, 0.02379263189520464, 0.7388214356402985, 0.01880438318540079),
    @("RTREE", 0.7760602454867427, 0.02453797399747067, 0.7698881594501135, 0.02458547155843596, 0.7514536085385257, 0.01357839011082922, 0.7502698101949177, 0.01085075024341467, 0.7534135534385253, 0.03401072561664792),
    @("XTREE", 0.8334858289152438, 0.01803767283039963, 0.820721388207543,  0.01724801362350073, 0.8122686183217327, 0.02342287613436309, 0.7954826323819981, 0.02826668479107708, 0.7936592651557134, 0.02099679302383254),
    @("SVM",   0.8343199710647928, 0.01102846655358631, 0.8214478402825881, 0.01406192486750847, 0.8143759897763273, 0.0220104682079383,  0.8091005566881602, 0.023838555299636,   0.797632064368036,  0.01156208009881626)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = 2 + $r
    $entry = $data[$r]

    # Column A: zero-based index
    $ws.Cells.Item($rowNum, 1).Value = $r

    # Column B: algorithm name
    $ws.Cells.Item($rowNum, 2).Value = $entry[0]

    # Columns C..L: the ten numeric values
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($rowNum, 3 + $c).Value = $entry[1 + $c]
    }
}

# Remove the now-unused row 9 (previously held the NB algorithm, which is removed)
$ws.Range("A9:L9").Delete()
